# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell (outside the used range) used to push values in as literal
# text via Copy/PasteSpecial -- this avoids Excel auto-converting
# numeric-looking strings (e.g. "608.58") into real numbers while leaving
# the destination cells existing style/format untouched.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

$ws.Range("D2").Value = "67.728.88"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "3.511.17"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  -0.03%  "
$scratch.Value = "608.58"
$scratch.Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$ws.Range("E5").Value = "  -1.08%  "
$scratch.Value = "152.23"
$scratch.Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "3.509.17"
$ws.Range("E7").Value = "  -0.52%  "
$ws.Range("E8").Value = "  +0.04%  "
$scratch.Value = "0.488"
$scratch.Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4163) | Out-Null
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("E10").Value = "  +2.86%  "
$scratch.Value = "7.64"
$scratch.Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4163) | Out-Null
$ws.Range("E11").Value = "  +7.43%  "
$scratch.Value = "0.433"
$scratch.Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4163) | Out-Null
$ws.Range("E12").Value = "  +1.58%  "
$scratch.Value = "0.0000218"
$scratch.Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4163) | Out-Null
$ws.Range("E13").Value = "  -1.84%  "
$scratch.Value = "32.27"
$scratch.Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4163) | Out-Null
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "4.101.47"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").Value = "67.625.41"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "3.483.50"
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("E18").Value = "  -0.67%  "
$scratch.Value = "6.53"
$scratch.Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4163) | Out-Null
$ws.Range("E19").Value = "  +1.74%  "
$ws.Range("E20").Value = "  +1.12%  "
$scratch.Value = "9.89"
$scratch.Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4163) | Out-Null
$ws.Range("E21").Value = "  +3.84%  "
$scratch.Value = "447.75"
$scratch.Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("E23").Value = "  +0.66%  "
$scratch.Value = "78.46"
$scratch.Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4163) | Out-Null
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").Value = "3.648.84"
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$scratch.Value = "0.0000128"
$scratch.Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4163) | Out-Null
$ws.Range("E26").Value = "  -3.29%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$scratch.Value = "1.00"
$scratch.Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$scratch.Value = "8.81"
$scratch.Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4163) | Out-Null
$ws.Range("E28").Value = "  +3.36%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$scratch.Value = "10.09"
$scratch.Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$ws.Range("E29").Value = "  -1.73%  "
$scratch.Value = "2.52"
$scratch.Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163) | Out-Null
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").Value = "  +5.24%  "
$ws.Range("E32").Value = "  +4.27%  "
$ws.Range("E33").Value = "  +0.09%  "
$scratch.Value = "25.62"
$scratch.Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4163) | Out-Null
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("E35").Value = "  -0.53%  "
$scratch.Value = "1.87"
$scratch.Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4163) | Out-Null
$ws.Range("E36").Value = "  +0.97%  "
$ws.Range("D37").Value = "3.500.28"
$ws.Range("E37").Value = "  -0.55%  "
$scratch.Value = "7.99"
$scratch.Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null
$ws.Range("E40").Value = "  +6.72%  "
$scratch.Value = "178.93"
$scratch.Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163) | Out-Null
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("E42").Value = "  +0.02%  "
$scratch.Value = "0.0897"
$scratch.Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4163) | Out-Null
$ws.Range("E43").Value = "  +1.46%  "
$ws.Range("E44").Value = "  +0.16%  "
$scratch.Value = "0.892"
$scratch.Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$ws.Range("E45").Value = "  +1.12%  "
$scratch.Value = "30.13"
$scratch.Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163) | Out-Null
$ws.Range("E46").Value = "  +5.60%  "
$scratch.Value = "46.48"
$scratch.Copy() | Out-Null
$ws.Range("D47").PasteSpecial(-4163) | Out-Null
$ws.Range("E47").Value = "  +2.95%  "
$scratch.Value = "1.31"
$scratch.Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163) | Out-Null
$ws.Range("E48").Value = "  +3.75%  "
$ws.Range("E49").Value = "  -2.60%  "
$scratch.Value = "7.63"
$scratch.Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4163) | Out-Null
$ws.Range("E50").Value = "  +0.26%  "
$scratch.Value = "0.254"
$scratch.Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4163) | Out-Null
$ws.Range("E51").Value = "  +1.96%  "

# Remove scratch cell content/format so it does not alter the used range
$scratch.Clear() | Out-Null
$excel.CutCopyMode = $false
